# SetCellValueTypes.xlsx sample refresh:
#  - drop the "Xceed Trial License" watermark sheet (trial banner removed)
#  - regenerate the sample's guid: value
#  - regenerate the sample's DateTime: value
#  - minor column width relayout on Sheet1 (column H)

$wb = $excel.ActiveWorkbook

# Deleting a sheet normally prompts a confirmation dialog - suppress it.
$excel.DisplayAlerts = $false

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Xceed Trial License")

# Remove the trial-license sheet entirely (and its content/relationships).
$ws2.Delete()

# Make Sheet1 the active/selected sheet and reset the view to the default.
$ws1.Select()
$ws1.Range("A1").Select()

# Refresh the sample's generated GUID value (H8, under the "guid:" label).
$ws1.Range("H8").Value = "3847d8cd-545a-4cdc-993c-ef6be2de72ef"

# Refresh the sample's generated DateTime value (E5, under "DateTime:").
$ws1.Range("E5").Value = 45567.560114537

# Small relayout of column H's width from the regenerated sample.
$ws1.Columns.Item(8).ColumnWidth = 37.1

$excel.DisplayAlerts = $true
